$wb = $excel.ActiveWorkbook

# --- Playoff Odds sheet: update Monte Carlo simulation results (fixed ties handling) ---
$wsPlayoff = $wb.Worksheets.Item("Playoff Odds")

# Rows 6 and 7 team order swapped (Pho_King and Eli exchange positions)
$wsPlayoff.Cells.Item(6, 1).Value = "Pho_King"
$wsPlayoff.Cells.Item(7, 1).Value = "Eli"

# Row 2
$wsPlayoff.Cells.Item(2, 2).Value = 41.2
$wsPlayoff.Cells.Item(2, 3).Value = 24
$wsPlayoff.Cells.Item(2, 4).Value = 13.5
$wsPlayoff.Cells.Item(2, 5).Value = 8.7
$wsPlayoff.Cells.Item(2, 6).Value = 5.8
$wsPlayoff.Cells.Item(2, 7).Value = 3.5
$wsPlayoff.Cells.Item(2, 8).Value = 1
$wsPlayoff.Cells.Item(2, 9).Value = 1.4
$wsPlayoff.Cells.Item(2, 10).Value = 0.6
$wsPlayoff.Cells.Item(2, 11).Value = 0.3
$wsPlayoff.Cells.Item(2, 12).Value = 0
$wsPlayoff.Cells.Item(2, 13).Value = 0
$wsPlayoff.Cells.Item(2, 14).Value = 96.7

# Row 3
$wsPlayoff.Cells.Item(3, 2).Value = 27.6
$wsPlayoff.Cells.Item(3, 3).Value = 20.3
$wsPlayoff.Cells.Item(3, 4).Value = 15
$wsPlayoff.Cells.Item(3, 5).Value = 10
$wsPlayoff.Cells.Item(3, 6).Value = 9.2
$wsPlayoff.Cells.Item(3, 7).Value = 5.4
$wsPlayoff.Cells.Item(3, 8).Value = 4.6
$wsPlayoff.Cells.Item(3, 9).Value = 3.4
$wsPlayoff.Cells.Item(3, 10).Value = 2.9
$wsPlayoff.Cells.Item(3, 11).Value = 1
$wsPlayoff.Cells.Item(3, 12).Value = 0.5
$wsPlayoff.Cells.Item(3, 13).Value = 0.1
$wsPlayoff.Cells.Item(3, 14).Value = 87.5

# Row 4
$wsPlayoff.Cells.Item(4, 2).Value = 10.9
$wsPlayoff.Cells.Item(4, 3).Value = 19
$wsPlayoff.Cells.Item(4, 4).Value = 18.9
$wsPlayoff.Cells.Item(4, 5).Value = 16.2
$wsPlayoff.Cells.Item(4, 6).Value = 10.6
$wsPlayoff.Cells.Item(4, 7).Value = 10.1
$wsPlayoff.Cells.Item(4, 8).Value = 5.7
$wsPlayoff.Cells.Item(4, 9).Value = 3
$wsPlayoff.Cells.Item(4, 10).Value = 3
$wsPlayoff.Cells.Item(4, 11).Value = 1.6
$wsPlayoff.Cells.Item(4, 12).Value = 0.9
$wsPlayoff.Cells.Item(4, 13).Value = 0.1
$wsPlayoff.Cells.Item(4, 14).Value = 85.7

# Row 5
$wsPlayoff.Cells.Item(5, 2).Value = 7.9
$wsPlayoff.Cells.Item(5, 3).Value = 9.7
$wsPlayoff.Cells.Item(5, 4).Value = 13.8
$wsPlayoff.Cells.Item(5, 5).Value = 12.6
$wsPlayoff.Cells.Item(5, 6).Value = 12.4
$wsPlayoff.Cells.Item(5, 7).Value = 12.4
$wsPlayoff.Cells.Item(5, 8).Value = 9.9
$wsPlayoff.Cells.Item(5, 9).Value = 7.7
$wsPlayoff.Cells.Item(5, 10).Value = 6
$wsPlayoff.Cells.Item(5, 11).Value = 4.3
$wsPlayoff.Cells.Item(5, 12).Value = 2.5
$wsPlayoff.Cells.Item(5, 13).Value = 0.8
$wsPlayoff.Cells.Item(5, 14).Value = 68.8

# Row 6
$wsPlayoff.Cells.Item(6, 2).Value = 4.1
$wsPlayoff.Cells.Item(6, 3).Value = 7.9
$wsPlayoff.Cells.Item(6, 4).Value = 9.1
$wsPlayoff.Cells.Item(6, 5).Value = 13
$wsPlayoff.Cells.Item(6, 6).Value = 11.2
$wsPlayoff.Cells.Item(6, 7).Value = 12.9
$wsPlayoff.Cells.Item(6, 8).Value = 11.8
$wsPlayoff.Cells.Item(6, 9).Value = 10.2
$wsPlayoff.Cells.Item(6, 10).Value = 8
$wsPlayoff.Cells.Item(6, 11).Value = 7
$wsPlayoff.Cells.Item(6, 12).Value = 2.9
$wsPlayoff.Cells.Item(6, 13).Value = 1.9
$wsPlayoff.Cells.Item(6, 14).Value = 58.2

# Row 7
$wsPlayoff.Cells.Item(7, 2).Value = 4.8
$wsPlayoff.Cells.Item(7, 3).Value = 7.4
$wsPlayoff.Cells.Item(7, 4).Value = 9.7
$wsPlayoff.Cells.Item(7, 5).Value = 10
$wsPlayoff.Cells.Item(7, 6).Value = 13.8
$wsPlayoff.Cells.Item(7, 7).Value = 12.5
$wsPlayoff.Cells.Item(7, 8).Value = 13.2
$wsPlayoff.Cells.Item(7, 9).Value = 10.4
$wsPlayoff.Cells.Item(7, 10).Value = 8.2
$wsPlayoff.Cells.Item(7, 11).Value = 5.8
$wsPlayoff.Cells.Item(7, 12).Value = 3.5
$wsPlayoff.Cells.Item(7, 13).Value = 0.7
$wsPlayoff.Cells.Item(7, 14).Value = 58.2

# Row 8
$wsPlayoff.Cells.Item(8, 2).Value = 1.4
$wsPlayoff.Cells.Item(8, 3).Value = 3.4
$wsPlayoff.Cells.Item(8, 4).Value = 7.4
$wsPlayoff.Cells.Item(8, 5).Value = 10
$wsPlayoff.Cells.Item(8, 6).Value = 10.7
$wsPlayoff.Cells.Item(8, 7).Value = 12.6
$wsPlayoff.Cells.Item(8, 8).Value = 13.8
$wsPlayoff.Cells.Item(8, 9).Value = 13
$wsPlayoff.Cells.Item(8, 10).Value = 12.6
$wsPlayoff.Cells.Item(8, 11).Value = 8.6
$wsPlayoff.Cells.Item(8, 12).Value = 4.9
$wsPlayoff.Cells.Item(8, 13).Value = 1.6
$wsPlayoff.Cells.Item(8, 14).Value = 45.5

# Row 9
$wsPlayoff.Cells.Item(9, 2).Value = 1.6
$wsPlayoff.Cells.Item(9, 3).Value = 4.6
$wsPlayoff.Cells.Item(9, 4).Value = 4.9
$wsPlayoff.Cells.Item(9, 5).Value = 8.4
$wsPlayoff.Cells.Item(9, 6).Value = 10.5
$wsPlayoff.Cells.Item(9, 7).Value = 10.8
$wsPlayoff.Cells.Item(9, 8).Value = 11.4
$wsPlayoff.Cells.Item(9, 9).Value = 12.3
$wsPlayoff.Cells.Item(9, 10).Value = 12.6
$wsPlayoff.Cells.Item(9, 11).Value = 10.5
$wsPlayoff.Cells.Item(9, 12).Value = 7.5
$wsPlayoff.Cells.Item(9, 13).Value = 4.9
$wsPlayoff.Cells.Item(9, 14).Value = 40.8

# Row 10
$wsPlayoff.Cells.Item(10, 2).Value = 0.3
$wsPlayoff.Cells.Item(10, 3).Value = 2.1
$wsPlayoff.Cells.Item(10, 4).Value = 4.8
$wsPlayoff.Cells.Item(10, 5).Value = 6.1
$wsPlayoff.Cells.Item(10, 6).Value = 9.2
$wsPlayoff.Cells.Item(10, 7).Value = 7.9
$wsPlayoff.Cells.Item(10, 8).Value = 11.5
$wsPlayoff.Cells.Item(10, 9).Value = 13.9
$wsPlayoff.Cells.Item(10, 10).Value = 13.2
$wsPlayoff.Cells.Item(10, 11).Value = 13.1
$wsPlayoff.Cells.Item(10, 12).Value = 11.2
$wsPlayoff.Cells.Item(10, 13).Value = 6.7
$wsPlayoff.Cells.Item(10, 14).Value = 30.4

# Row 11
$wsPlayoff.Cells.Item(11, 2).Value = 0.1
$wsPlayoff.Cells.Item(11, 3).Value = 0.8
$wsPlayoff.Cells.Item(11, 4).Value = 1.6
$wsPlayoff.Cells.Item(11, 5).Value = 2.2
$wsPlayoff.Cells.Item(11, 6).Value = 3
$wsPlayoff.Cells.Item(11, 7).Value = 4.9
$wsPlayoff.Cells.Item(11, 8).Value = 6.4
$wsPlayoff.Cells.Item(11, 9).Value = 10
$wsPlayoff.Cells.Item(11, 10).Value = 12.4
$wsPlayoff.Cells.Item(11, 11).Value = 17
$wsPlayoff.Cells.Item(11, 12).Value = 22.5
$wsPlayoff.Cells.Item(11, 13).Value = 19.1
$wsPlayoff.Cells.Item(11, 14).Value = 12.6

# Row 12
$wsPlayoff.Cells.Item(12, 2).Value = 0.1
$wsPlayoff.Cells.Item(12, 3).Value = 0.7
$wsPlayoff.Cells.Item(12, 4).Value = 1
$wsPlayoff.Cells.Item(12, 5).Value = 2.2
$wsPlayoff.Cells.Item(12, 6).Value = 1.9
$wsPlayoff.Cells.Item(12, 7).Value = 5.1
$wsPlayoff.Cells.Item(12, 8).Value = 6.3
$wsPlayoff.Cells.Item(12, 9).Value = 8.3
$wsPlayoff.Cells.Item(12, 10).Value = 10.9
$wsPlayoff.Cells.Item(12, 11).Value = 18.1
$wsPlayoff.Cells.Item(12, 12).Value = 19.7
$wsPlayoff.Cells.Item(12, 13).Value = 25.7
$wsPlayoff.Cells.Item(12, 14).Value = 11

# Row 13
$wsPlayoff.Cells.Item(13, 2).Value = 0
$wsPlayoff.Cells.Item(13, 3).Value = 0.1
$wsPlayoff.Cells.Item(13, 4).Value = 0.3
$wsPlayoff.Cells.Item(13, 5).Value = 0.6
$wsPlayoff.Cells.Item(13, 6).Value = 1.7
$wsPlayoff.Cells.Item(13, 7).Value = 1.9
$wsPlayoff.Cells.Item(13, 8).Value = 4.4
$wsPlayoff.Cells.Item(13, 9).Value = 6.4
$wsPlayoff.Cells.Item(13, 10).Value = 9.6
$wsPlayoff.Cells.Item(13, 11).Value = 12.7
$wsPlayoff.Cells.Item(13, 12).Value = 23.9
$wsPlayoff.Cells.Item(13, 13).Value = 38.4
$wsPlayoff.Cells.Item(13, 14).Value = 4.6

# --- Record Odds sheet: update playoff chance, expected/most-likely records ---
$wsRecord = $wb.Worksheets.Item("Record Odds")

# Row 2
$wsRecord.Cells.Item(2, 6).Value = 97.7
$wsRecord.Cells.Item(2, 7).Value = "9.5-4.4-0.1"
$wsRecord.Cells.Item(2, 8).Value = "10-4"

# Row 3
$wsRecord.Cells.Item(3, 6).Value = 92.10000000000001
$wsRecord.Cells.Item(3, 7).Value = "8.9-5.1-0.1"
$wsRecord.Cells.Item(3, 8).Value = "9-5"

# Row 4
$wsRecord.Cells.Item(4, 6).Value = 91.4
$wsRecord.Cells.Item(4, 7).Value = "8.3-5.6-0.1"
$wsRecord.Cells.Item(4, 8).Value = "8-6"

# Row 5
$wsRecord.Cells.Item(5, 6).Value = 78.7
$wsRecord.Cells.Item(5, 7).Value = "7.6-6.3-0.1"
$wsRecord.Cells.Item(5, 8).Value = "8-6"

# Row 6
$wsRecord.Cells.Item(6, 6).Value = 71.39999999999999
$wsRecord.Cells.Item(6, 7).Value = "7.5-6.4-0.1"
$wsRecord.Cells.Item(6, 8).Value = "7-7"

# Row 7
$wsRecord.Cells.Item(7, 6).Value = 70
$wsRecord.Cells.Item(7, 7).Value = "7.5-6.4-0.1"
$wsRecord.Cells.Item(7, 8).Value = "8-6"

# Row 8
$wsRecord.Cells.Item(8, 6).Value = 59.3
$wsRecord.Cells.Item(8, 7).Value = "6.8-7.1-0.1"
$wsRecord.Cells.Item(8, 8).Value = "7-7"

# Row 9
$wsRecord.Cells.Item(9, 6).Value = 52.2
$wsRecord.Cells.Item(9, 7).Value = "6.3-7.6-0.1"
$wsRecord.Cells.Item(9, 8).Value = "6-8"

# Row 10
$wsRecord.Cells.Item(10, 6).Value = 41.9
$wsRecord.Cells.Item(10, 7).Value = "6.3-7.6-0.1"
$wsRecord.Cells.Item(10, 8).Value = "6-8"

# Row 11
$wsRecord.Cells.Item(11, 6).Value = 19
$wsRecord.Cells.Item(11, 7).Value = "5.2-8.8-0.1"
$wsRecord.Cells.Item(11, 8).Value = "5-9"

# Row 12
$wsRecord.Cells.Item(12, 6).Value = 17.3
$wsRecord.Cells.Item(12, 7).Value = "5.1-8.8-0.1"
$wsRecord.Cells.Item(12, 8).Value = "5-9"

# Row 13
$wsRecord.Cells.Item(13, 6).Value = 9
$wsRecord.Cells.Item(13, 7).Value = "4.6-9.4-0.1"
$wsRecord.Cells.Item(13, 8).Value = "5-9"
